# Added scanning page. Updated total ugc count to use SteamDB top 100 games
#
# Appends the "scanning page" worth of rows (14-25) which duplicate the
# existing Name/Creator/Country/Language/TUS/Comments/Rating columns for
# rows 2-13, but record the "Date Posted" column (H) as a real Excel date
# serial value (formatted mm-dd-yy) instead of free-text like "09/02/2024".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New rows 14-25: same level data as rows 2-13, with column H now a real
# date value (so it can be used for the SteamDB-top-100 re-scan ordering)
# ---------------------------------------------------------------------

$ws.Range("A14").Value = "HAPPY LAND:The Spring Festival"
$ws.Range("B14").Value = "激辣浪味仙☆"
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "en"
$ws.Range("E14").Value = 4644
$ws.Range("F14").Value = "N/A (Needs log in)"
$ws.Range("G14").Value = "7"
$ws.Range("H14").Value = 45300
$ws.Range("H14").NumberFormat = "mm-dd-yy"

$ws.Range("A15").Value = "墨墨28（龙年大吉）"
$ws.Range("B15").Value = "墨墨"
$ws.Range("C15").Value = "N/A"
$ws.Range("D15").Value = "zh-cn"
$ws.Range("E15").Value = 1284
$ws.Range("F15").Value = "N/A (Needs log in)"
$ws.Range("G15").Value = "4"
$ws.Range("H15").Value = 45299
$ws.Range("H15").NumberFormat = "mm-dd-yy"

$ws.Range("A16").Value = "2024happy new year"
$ws.Range("B16").Value = "脸红"
$ws.Range("C16").Value = "CN"
$ws.Range("D16").Value = "zh-cn"
$ws.Range("E16").Value = 2609
$ws.Range("F16").Value = "N/A (Needs log in)"
$ws.Range("G16").Value = "2"
$ws.Range("H16").Value = 45296
$ws.Range("H16").NumberFormat = "mm-dd-yy"

$ws.Range("A17").Value = "time tunnel"
$ws.Range("B17").Value = "MoQiangShou"
$ws.Range("C17").Value = "CN"
$ws.Range("D17").Value = "zh-cn"
$ws.Range("E17").Value = 3610
$ws.Range("F17").Value = "N/A (Needs log in)"
$ws.Range("G17").Value = "2"
$ws.Range("H17").Value = 45295
$ws.Range("H17").NumberFormat = "mm-dd-yy"

$ws.Range("A18").Value = "陈怼怼的生日派对！"
$ws.Range("B18").Value = "可爱屁の梓逸ღ⊰"
$ws.Range("C18").Value = "N/A"
$ws.Range("D18").Value = "zh-cn"
$ws.Range("E18").Value = 564
$ws.Range("F18").Value = "N/A (Needs log in)"
$ws.Range("G18").Value = "0"
$ws.Range("H18").Value = 45305
$ws.Range("H18").NumberFormat = "mm-dd-yy"

$ws.Range("A19").Value = "龙年大吉"
$ws.Range("B19").Value = "Cikey"
$ws.Range("C19").Value = "CN"
$ws.Range("D19").Value = "zh-cn"
$ws.Range("E19").Value = 677
$ws.Range("F19").Value = "N/A (Needs log in)"
$ws.Range("G19").Value = "0"
$ws.Range("H19").Value = 45301
$ws.Range("H19").NumberFormat = "mm-dd-yy"

$ws.Range("A20").Value = "1x1=1"
$ws.Range("B20").Value = "TaiLmEaT"
$ws.Range("C20").Value = "KR"
$ws.Range("D20").Value = "so"
$ws.Range("E20").Value = 175
$ws.Range("F20").Value = "N/A (Needs log in)"
$ws.Range("G20").Value = "4"
$ws.Range("H20").Value = 45309
$ws.Range("H20").NumberFormat = "mm-dd-yy"

$ws.Range("A21").Value = "Welcome to Hell"
$ws.Range("B21").Value = "ΜOΛΩΝ ΛΑΒΕ"
$ws.Range("C21").Value = "CA"
$ws.Range("D21").Value = "en"
$ws.Range("E21").Value = 1288
$ws.Range("F21").Value = "N/A (Needs log in)"
$ws.Range("G21").Value = "1"
$ws.Range("H21").Value = 45300
$ws.Range("H21").NumberFormat = "mm-dd-yy"

$ws.Range("A22").Value = "연습맵2"
$ws.Range("B22").Value = "qustjdbs"
$ws.Range("C22").Value = "N/A"
$ws.Range("D22").Value = "ca"
$ws.Range("E22").Value = 27
$ws.Range("F22").Value = "N/A (Needs log in)"
$ws.Range("G22").Value = "0"
$ws.Range("H22").Value = 45311
$ws.Range("H22").NumberFormat = "mm-dd-yy"

$ws.Range("A23").Value = "小点点的生日连跳图"
$ws.Range("B23").Value = "35`""
$ws.Range("C23").Value = "N/A"
$ws.Range("D23").Value = "zh-cn"
$ws.Range("E23").Value = 92
$ws.Range("F23").Value = "N/A (Needs log in)"
$ws.Range("G23").Value = "0"
$ws.Range("H23").Value = 45308
$ws.Range("H23").NumberFormat = "mm-dd-yy"

$ws.Range("A24").Value = "大年初一快乐"
$ws.Range("B24").Value = "甘九"
$ws.Range("C24").Value = "N/A"
$ws.Range("D24").Value = "zh-cn"
$ws.Range("E24").Value = 418
$ws.Range("F24").Value = "N/A (Needs log in)"
$ws.Range("G24").Value = "0"
$ws.Range("H24").Value = 45302
$ws.Range("H24").NumberFormat = "mm-dd-yy"

$ws.Range("A25").Value = "Quiz Scene"
$ws.Range("B25").Value = "dbabicius"
$ws.Range("C25").Value = "LT"
$ws.Range("D25").Value = "en"
$ws.Range("E25").Value = 76
$ws.Range("F25").Value = "N/A (Needs log in)"
$ws.Range("G25").Value = "0"
$ws.Range("H25").Value = 45318
$ws.Range("H25").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------
# Column widths: best-fit widths picked up by Excel for the columns that
# now hold the widest content (level name, creator, "N/A (Needs log in)",
# rating, the new mm-dd-yy date and the "Creator we are planning to
# approach" note column).
# ---------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 29.17
$ws.Columns.Item(2).ColumnWidth = 15.83
$ws.Columns.Item(6).ColumnWidth = 16.83
$ws.Columns.Item(7).ColumnWidth = 9.67
$ws.Columns.Item(8).ColumnWidth = 10.83
$ws.Columns.Item(10).ColumnWidth = 33

# ---------------------------------------------------------------------
# Scroll the view down to the newly-added data and leave the selection
# on the next empty date cell, like the author did after typing the rows.
# ---------------------------------------------------------------------

$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("H26").Select()
